$d = $word.ActiveDocument

function Replace-FirstMatch($searchText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if ($found) {
        $rng.Text = $replaceText
    }
    return $found
}

# 1. Order number on cover: No.  135 -> No.  137
Replace-FirstMatch "No.  135" "No.  137"

# 2. Both "Fecha:" and "Fecha evento:" dates: 2019-04-22 -> 2019-05-09 (two occurrences)
Replace-FirstMatch "2019-04-22" "2019-05-09"
Replace-FirstMatch "2019-04-22" "2019-05-09"

# 3. Damage description line: update invoice number and description text
Replace-FirstMatch "Pago por daños de compromiso con factura No 135: Pantaló dañado" "Pago por daños de compromiso con factura No 137: algo"

# 4. Both VALOR DAÑO amounts: 30000 -> 29999 (two occurrences)
Replace-FirstMatch "30000" "29999"
Replace-FirstMatch "30000" "29999"

# 5. Fecha devolución: 2019-04-23 -> 2019-05-10
Replace-FirstMatch "2019-04-23" "2019-05-10"
